$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The analyte list on Sheet1 runs down column A (A1 is the "Analyte" header,
# A2:A15 are the existing analytes). Add one more analyte, "test", in the
# next empty row, using the same look (font/alignment/row height) as the
# rest of the list - this is used to warn when an analyte from this list
# is missing from the data.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp
$newRow = $lastRow + 1

# Copy the formatting (style) of the last data row onto the new row.
$ws.Range("A" + $lastRow).Copy()
$ws.Range("A" + $newRow).PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("A" + $newRow).Value = "test"

# Match the row height that Excel auto-applies for this style (12pt font).
$ws.Rows.Item($newRow).RowHeight = $ws.Rows.Item($lastRow).RowHeight

$ws.Range("A" + $newRow).Select()
